# "Search test cases modified" - add three new test case rows (TestCase_A21..A23)
# to the "Test Cases" sheet, matching formatting used by neighbouring rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Copy cell formatting from existing rows into the new rows (22-24) ---
# Row 22 formatting: A<-A2 (s=2), B<-B6 (s=6), C<-C12 (s=9-equivalent), D<-D6 (s=6), E<-E2 (s=2)
$ws.Range("A2").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# Row 23 formatting: A<-A6 (s=6), B<-B6 (s=6), C<-C12 (s=9-equivalent), D<-D6 (s=6), E<-B6-like (s=6)
$ws.Range("A6").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("E23").PasteSpecial(-4122)

# Row 24 formatting: A<-A6 (s=6), B<-B6 (s=6), C<-A2-like (s=2), D<-D6 (s=6), E<-E2 (s=2)
$ws.Range("A6").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E24").PasteSpecial(-4122)

# --- Fill in the new values (order matches how the shared-string table grows) ---
$ws.Range("A22").Value = "TestCase_A21"
$ws.Range("C22").Value = "Verify View additional email preferences link is working"
$ws.Range("A23").Value = "TestCase_A22"
$ws.Range("A24").Value = "TestCase_A23"
$ws.Range("C24").Value = "Verify change password link in the account page is working correctly."
$ws.Range("C23").Value = 'Verify that the  checkbox  is present and can be modified for "Receive email notifications for likes,comments and other activity" is working correctly.'
$ws.Range("B23").Value = "OPQA-854,OPQA-853"
$ws.Range("B22").Value = "OPQA-399"
$ws.Range("B24").Value = "OPQA-527"

$ws.Range("D22").Value = "Y"
$ws.Range("E22").Value = "SKIP"
$ws.Range("D23").Value = "Y"
$ws.Range("E23").Value = "SKIP"
$ws.Range("D24").Value = "Y"
$ws.Range("E24").Value = "SKIP"

# Row 23 needs extra height to fit the wrapped text, like row 20 does.
$ws.Rows.Item(23).RowHeight = 30

# --- Refresh column B's auto-fit width now that it holds longer Jira ids, ---
# --- and column E's width, to mirror Excel's own recalculation. ---
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(5).AutoFit()

# --- Leftover view state from the last user interaction while editing ---
$ws.Range("D2:D24").Select()
$excel.ActiveWindow.ScrollColumn = 2

# --- Page setup (portrait print orientation), mirroring other sheets ---
$ws.PageSetup.Orientation = 1
